# Update workbook "上海-漫展信息.xlsx" to match the upstream data refresh.
# Sheets: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local Life),
#         4=全部类型 (All Types, a merged/aggregated view).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------------
# Sheet 2 (演出): two brand new events were scraped (2024-10-11 and
# 2024-10-13), inserted right before the existing 2024-10-18 "ROOKiEZ is
# PUNK`D" row, pushing the old rows 25-34 down to 27-36.
# ---------------------------------------------------------------------------
$ws2.Range("A25:A26").EntireRow.Insert()

# The insert clones formatting from the row below, which synthesizes a new
# (near-duplicate) style for column A. Re-paste the canonical bold/centered
# style used by every other row's index cell so the new cells reuse the
# original style id instead of a duplicate one.
$ws2.Range("A1").Copy()
$ws2.Range("A25:A26").PasteSpecial(-4122)

# Column B stores plain text dates (e.g. "2024-10-11"); without an explicit
# text format Excel's COM layer auto-converts such literals into date
# serials, so force the text format first.
$ws2.Range("B25:B36").NumberFormat = "@"

# New row 25: 上海·井草圣二 2024《落叶轻扬》指弹吉他音乐会
$ws2.Range("A25").Value = 24
$ws2.Range("B25").Value = "2024-10-11"
$ws2.Range("C25").Value = "上海·井草圣二 2024《落叶轻扬》指弹吉他音乐会"
$ws2.Range("D25").Value = "宜昌路179号2F 万代南梦宫上海文化中心-未来剧场"
$ws2.Range("E25").Value = "2024.10.11 19:30-10.11 21:30"
$ws2.Range("F25").Value = 0
$ws2.Range("G25").Value = 260
$ws2.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=91647"
$ws2.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202409/Y9IYLEwL1725247462745.jpeg"

# New row 26: 上海·Luca Stricagnoli 2024《进化时间》指弹吉他音乐会
$ws2.Range("A26").Value = 25
$ws2.Range("B26").Value = "2024-10-13"
$ws2.Range("C26").Value = "上海·Luca Stricagnoli 2024《进化时间》指弹吉他音乐会"
$ws2.Range("D26").Value = "宜昌路179号2F 万代南梦宫上海文化中心-未来剧场"
$ws2.Range("E26").Value = "2024.10.13 19:30-10.13 21:30"
$ws2.Range("F26").Value = 0
$ws2.Range("G26").Value = 220
$ws2.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=91645"
$ws2.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202409/tCnThnQ41725246720534.png"

# The shifted rows (old 25-34, now 27-36) keep all of their original data via
# the row insert/shift; only the manually-maintained running index in column
# A needs to be bumped by 2 to stay sequential (A = row - 1).
$ws2.Range("A27").Value = 26
$ws2.Range("A28").Value = 27
$ws2.Range("A29").Value = 28
$ws2.Range("A30").Value = 29
$ws2.Range("A31").Value = 30
$ws2.Range("A32").Value = 31
$ws2.Range("A33").Value = 32
$ws2.Range("A34").Value = 33
$ws2.Range("A35").Value = 34
$ws2.Range("A36").Value = 35

# ---------------------------------------------------------------------------
# Sheet 1 (展览): refreshed view/interest counters.
# ---------------------------------------------------------------------------
$ws1.Range("F3").Value = 3295
$ws1.Range("F6").Value = 7703
$ws1.Range("F8").Value = 696
$ws1.Range("F9").Value = 1123
$ws1.Range("F13").Value = 1715
$ws1.Range("G13").Value = 98
$ws1.Range("F15").Value = 6140
$ws1.Range("F16").Value = 63
$ws1.Range("F19").Value = 1021
$ws1.Range("F20").Value = 1015
$ws1.Range("F21").Value = 4256
$ws1.Range("F22").Value = 5399
$ws1.Range("F23").Value = 340
$ws1.Range("F24").Value = 150
$ws1.Range("F25").Value = 1067
$ws1.Range("F32").Value = 95
$ws1.Range("F34").Value = 427
$ws1.Range("F37").Value = 55
$ws1.Range("F38").Value = 594
$ws1.Range("F39").Value = 407
$ws1.Range("F40").Value = 330
$ws1.Range("F41").Value = 1168
$ws1.Range("F44").Value = 3177
$ws1.Range("F47").Value = 37

# ---------------------------------------------------------------------------
# Sheet 2 (演出): refreshed view/interest counters on the untouched rows.
# ---------------------------------------------------------------------------
$ws2.Range("F5").Value = 380
$ws2.Range("F6").Value = 633
$ws2.Range("F11").Value = 268
$ws2.Range("F16").Value = 169
$ws2.Range("F24").Value = 6499

# ---------------------------------------------------------------------------
# Sheet 3 (本地生活): refreshed view/interest counters, plus one status
# string flip from "已售罄" (sold out) to "暂时售罄" (temporarily sold out).
# ---------------------------------------------------------------------------
$ws3.Range("F6").Value = 1290
$ws3.Range("F8").Value = 542
$ws3.Range("F9").Value = 2108
$ws3.Range("G9").Value = "暂时售罄"
$ws3.Range("F10").Value = 8852
$ws3.Range("F11").Value = 956
$ws3.Range("F12").Value = 70

# ---------------------------------------------------------------------------
# Sheet 4 (全部类型): same refreshed counters mirrored into the merged view
# (this sheet is not affected by the two new rows inserted in sheet 2).
# ---------------------------------------------------------------------------
$ws4.Range("F3").Value = 3295
$ws4.Range("F6").Value = 542
$ws4.Range("F7").Value = 2108
$ws4.Range("G7").Value = "暂时售罄"
$ws4.Range("F8").Value = 956
$ws4.Range("F9").Value = 380
$ws4.Range("F10").Value = 70
$ws4.Range("F11").Value = 696
$ws4.Range("F12").Value = 1123
$ws4.Range("F18").Value = 6140
$ws4.Range("F19").Value = 63
$ws4.Range("F21").Value = 1021
$ws4.Range("F22").Value = 1015
$ws4.Range("F23").Value = 4256
$ws4.Range("F24").Value = 5399
$ws4.Range("F25").Value = 340
$ws4.Range("F26").Value = 150
$ws4.Range("F27").Value = 1067
$ws4.Range("F31").Value = 95
$ws4.Range("F33").Value = 427
$ws4.Range("F35").Value = 169
$ws4.Range("F37").Value = 55
$ws4.Range("F38").Value = 594
$ws4.Range("F39").Value = 407
$ws4.Range("F40").Value = 330
$ws4.Range("F43").Value = 3177
$ws4.Range("F46").Value = 38
